$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column L: quantity-per-32 factor (header constant) and per-row formula
$ws.Range("L1").Value = 32
$ws.Range("L2").Formula = "=F2*L$1"
$ws.Range("L3:L13").Formula = "=F3*L$1"

# New column M: price/value figures per row (no formula, literal numbers)
$ws.Range("M2").Value = 800
$ws.Range("M3").Value = 74
$ws.Range("M4").Value = 30
$ws.Range("M5").Value = 32
$ws.Range("M6").Value = 30
$ws.Range("M7").Value = 0
$ws.Range("M8").Value = 1050
$ws.Range("M9").Value = 10
$ws.Range("M10").Value = 128
$ws.Range("M12").Value = 64
$ws.Range("M13").Value = 29

# New column N: a couple of extra figures
$ws.Range("N9").Value = 100
$ws.Range("N13").Value = 1

# Restore selection/active cell as left by the author
$ws.Range("J13").Select() | Out-Null
